$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.183.75"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +0.56%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.685.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.15%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  +0.16%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'215.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = "'  +0.48%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("E7").Value = "'  +0.14%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'23.18"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +8.45%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.261"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +3.50%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +0.81%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +0.20%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'1.923.72"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.17%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.692.72"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +0.48%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'4.20"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +2.09%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.555"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +3.75%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'67.00"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +1.30%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'27.183.59"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +0.56%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'236.04"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.23%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'8.00"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -2.08%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.0₃0744"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +1.08%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("E22").Value = "'  +1.94%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'9.61"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +3.83%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -2.78%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'146.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +0.09%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'7.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.15%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'16.44"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +2.16%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  +0.28%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  +0.01%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +0.70%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -0.04%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  +1.33%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("E33").Value = "'  +1.52%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +1.83%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "'  -1.54%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("E36").Value = "'  +2.80%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.948"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  +2.82%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -0.30%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -0.69%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "'  +1.37%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'69.10"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  +1.13%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'  -0.23%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  +0.10%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  -0.88%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.833.09"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.49%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.790"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  +1.03%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'90.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -0.07%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +5.24%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("E49").Value = "'  +5.89%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'8.30"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +5.81%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = "'  -0.44%  "
$ws.Range("E51").Style = "Normal"
